# "trying out new math" - rework the Sheet1 math tables:
#  - rename headers, convert the x100/x20/... labels into plain multiplier numbers
#  - update the Coin/Multi count columns with the new distribution
#  - change the Total-weight formula range and relocate the grand-total cell
#  - add "Average multi value" / "Average coin value" boxed summary cells
#  - add a new block of named probability constants below the table

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 1 headers
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Coin"
$ws.Range("D1").Value = "Total weight"
$ws.Range("E1").Value = "Multi"
$ws.Range("H1").Value = "Total weight"

# Move "Total of those two" from J1/J2 over to M1/M2
$ws.Range("J1").ClearContents()
$ws.Range("J2").ClearContents()
$ws.Range("M1").Value = "Total of those two"

# ---------------------------------------------------------------------------
# Coin table (A:D) - counts in column B change for several rows
# ---------------------------------------------------------------------------
$ws.Range("B6").Value = 16000
$ws.Range("B7").Value = 39900
$ws.Range("B8").Value = 85000
$ws.Range("B9").Value = 130000
$ws.Range("B10").Value = 325000
$ws.Range("B11").Value = 475000
$ws.Range("B12").Value = 1450000
$ws.Range("B13").Value = 3000000
$ws.Range("B14").Value = 4306250

# ---------------------------------------------------------------------------
# Multi table (E:H) - E2:E7 switch from text labels ("x100" ...) to plain
# numbers, F column counts change, and the total-weight formula narrows
# from the whole column to F1:F10.
# ---------------------------------------------------------------------------
$ws.Range("E2").Value = 100
$ws.Range("E3").Value = 20
$ws.Range("E4").Value = 10
$ws.Range("E5").Value = 5
$ws.Range("E6").Value = 3
$ws.Range("E7").Value = 2

$ws.Range("F2").Value = 10
$ws.Range("F3").Value = 115
$ws.Range("F4").Value = 1200
$ws.Range("F5").Value = 9225
$ws.Range("F6").Value = 36450
$ws.Range("F7").Value = 123000

$ws.Range("H2").Formula = "=SUM(F1:F10)"

# Relocated grand-total cell
$ws.Range("M2").Formula = "=D2+H2"

# ---------------------------------------------------------------------------
# "Average multi value" boxed summary (row 8, columns E:G)
# ---------------------------------------------------------------------------
$avgMulti = $ws.Range("E8:G8")
$avgMulti.Font.Bold = $true
$avgMulti.Borders.Item(7).Weight = -4138    # xlEdgeLeft, xlMedium
$avgMulti.Borders.Item(8).Weight = -4138    # xlEdgeTop, xlMedium
$avgMulti.Borders.Item(9).Weight = -4138    # xlEdgeBottom, xlMedium
$avgMulti.Borders.Item(10).Weight = -4138   # xlEdgeRight, xlMedium

$ws.Range("E8").Value = "Average multi value"
$ws.Range("F8").ClearContents()
$ws.Range("G8").Formula = "=SUMPRODUCT(E2:E14, F2:F14)/SUM(F2:F14)"

# Drop the leftover "0.00%" placeholder styling on G10:G13 (kept on G9/G14)
$ws.Range("G10").Clear()
$ws.Range("G11").Clear()
$ws.Range("G12").Clear()
$ws.Range("G13").Clear()

# ---------------------------------------------------------------------------
# "Average coin value" boxed summary (new row 15, columns A:C)
# ---------------------------------------------------------------------------
$avgCoin = $ws.Range("A15:C15")
$avgCoin.Font.Bold = $true
$avgCoin.Borders.Item(7).Weight = -4138
$avgCoin.Borders.Item(8).Weight = -4138
$avgCoin.Borders.Item(9).Weight = -4138
$avgCoin.Borders.Item(10).Weight = -4138

$ws.Range("A15").Value = "Average coin value"
$ws.Range("C15").Formula = "=SUMPRODUCT(A2:A14, B2:B14)/SUM(B2:B14)"
$ws.Range("C15").NumberFormat = "0.000000"

# ---------------------------------------------------------------------------
# New block of named probability constants (rows 17-20)
# ---------------------------------------------------------------------------
$ws.Range("A17").Value = "Chance of getting C on every block (default value)"
$ws.Range("C17").NumberFormat = "0.00%"
$ws.Range("F17").Value = 0.0135
$ws.Range("F17").Font.Bold = $true

$ws.Range("A18").Value = "Chance of getting C decrease for every previous C"
$ws.Range("C18").NumberFormat = "0.00%"
$ws.Range("F18").Value = 0.425

$ws.Range("A19").Value = "Minimum value of Chance of getting C"
$ws.Range("C19").NumberFormat = "0.00%"
$ws.Range("F19").Value = 0.0015

$ws.Range("A20").Value = "(Multiply those numbers by 100 for % value)"
$ws.Range("C20").NumberFormat = "0.00%"

# ---------------------------------------------------------------------------
# Column width - column A now matches column B's 9.5 "bestFit" width
# ---------------------------------------------------------------------------
$ws.Range("A1").ColumnWidth = $ws.Range("B1").ColumnWidth

# ---------------------------------------------------------------------------
# Selection housekeeping
# ---------------------------------------------------------------------------
$ws.Range("F3").Select()
